# Natmi following Dr Hou advice
# Update the Gdf15-Gfral LR-pair results sheet (natmiOut/YoungD4/LR-pairs_lrc2p/Gdf15-Gfral.xlsx):
#  - recompute row 2 (ECs sending cluster) with new values
#  - insert a new "FAPs" sending-cluster row as row 3, shifting the former
#    M1/M2/sCs rows down to rows 4/5/6 with freshly recomputed values
#  - add a brand new "sCs" sending-cluster row as row 6

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> Gdf15 -> Gfral -> FAPs
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Gdf15"
$ws.Range("C2").Value = "Gfral"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.595593
$ws.Range("H2").Value = 4.786779
$ws.Range("I2").Value = 0.04211707711903992
$ws.Range("J2").Value = 0.04211707711903992
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.6050859999999999
$ws.Range("N2").Value = 1.815258
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 0.9654709859979999
$ws.Range("R2").Value = 8.689238873981999
$ws.Range("S2").Value = 0.04211707711903992
$ws.Range("T2").Value = 0.04211707711903992

# Row 3: FAPs -> Gdf15 -> Gfral -> FAPs  (new row, inserted before old M1 row)
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Gdf15"
$ws.Range("C3").Value = "Gfral"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.5295143333333333
$ws.Range("H3").Value = 1.588543
$ws.Range("I3").Value = 0.0139769953945881
$ws.Range("J3").Value = 0.0139769953945881
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.6050859999999999
$ws.Range("N3").Value = 1.815258
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 1
$ws.Range("Q3").Value = 0.3204017098993333
$ws.Range("R3").Value = 2.883615389094
$ws.Range("S3").Value = 0.0139769953945881
$ws.Range("T3").Value = 0.0139769953945881

# Row 4: M1 -> Gdf15 -> Gfral -> FAPs
$ws.Range("A4").Value = "M1"
$ws.Range("B4").Value = "Gdf15"
$ws.Range("C4").Value = "Gfral"
$ws.Range("D4").Value = "FAPs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 15.63002466666667
$ws.Range("H4").Value = 46.890074
$ws.Range("I4").Value = 0.4125682139859577
$ws.Range("J4").Value = 0.4125682139859577
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.6050859999999999
$ws.Range("N4").Value = 1.815258
$ws.Range("O4").Value = 1
$ws.Range("P4").Value = 1
$ws.Range("Q4").Value = 9.457509105454664
$ws.Range("R4").Value = 85.11758194909198
$ws.Range("S4").Value = 0.4125682139859577
$ws.Range("T4").Value = 0.4125682139859577

# Row 5: M2 -> Gdf15 -> Gfral -> FAPs
$ws.Range("A5").Value = "M2"
$ws.Range("B5").Value = "Gdf15"
$ws.Range("C5").Value = "Gfral"
$ws.Range("D5").Value = "FAPs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 19.20541533333333
$ws.Range("H5").Value = 57.616246
$ws.Range("I5").Value = 0.506943787480386
$ws.Range("J5").Value = 0.506943787480386
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.6050859999999999
$ws.Range("N5").Value = 1.815258
$ws.Range("O5").Value = 1
$ws.Range("P5").Value = 1
$ws.Range("Q5").Value = 11.62092794238533
$ws.Range("R5").Value = 104.588351481468
$ws.Range("S5").Value = 0.506943787480386
$ws.Range("T5").Value = 0.506943787480386

# Row 6: sCs -> Gdf15 -> Gfral -> FAPs  (brand new row)
$ws.Range("A6").Value = "sCs"
$ws.Range("B6").Value = "Gdf15"
$ws.Range("C6").Value = "Gfral"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.9241566666666666
$ws.Range("H6").Value = 2.77247
$ws.Range("I6").Value = 0.0243939260200282
$ws.Range("J6").Value = 0.0243939260200282
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.6050859999999999
$ws.Range("N6").Value = 1.815258
$ws.Range("O6").Value = 1
$ws.Range("P6").Value = 1
$ws.Range("Q6").Value = 0.5591942608066666
$ws.Range("R6").Value = 5.032748347259999
$ws.Range("S6").Value = 0.0243939260200282
$ws.Range("T6").Value = 0.0243939260200282
